$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("payment_methods")
$ws.Activate()
